$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 33 (Employment % of total): Micro / SMEs / MSMEs
$ws.Cells.Item(33, 2).Value = "'73.78"
$ws.Cells.Item(33, 3).Value = "'2.87"
$ws.Cells.Item(33, 4).Value = "'76.65"

# Row 34 (Employment absolute #): Micro / MSMEs (SMEs column 33.2 unchanged)
$ws.Cells.Item(34, 2).Value = "'38.48"
$ws.Cells.Item(34, 4).Value = "'71.68"

# Row 36 (Enterprises % of total): Micro / SMEs / MSMEs
$ws.Cells.Item(36, 2).Value = "'96.15"
$ws.Cells.Item(36, 3).Value = "'3.73"
$ws.Cells.Item(36, 4).Value = "'99.88"

# Row 40 (Value added to the economy % of total): Micro / SMEs / MSMEs
$ws.Cells.Item(40, 2).Value = "'26.44"
$ws.Cells.Item(40, 3).Value = "'37.63"
$ws.Cells.Item(40, 4).Value = "'64.07"
